$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 (2-RAP / SQUIMAN) is being removed; this deletes the whole row
# and shifts every subsequent row (16-19) up by one, which reproduces
# the remaining diffs to rows 15-18 and removes the old row 19 entirely.
$ws.Rows.Item(15).Delete()
